# Auto-generated edit script: update market-price-derived columns (H-N)
# across ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets per scheduled-runner refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(32, 8).Value = 2535.8572
$ws.Cells.Item(32, 9).Value = 1846.8334
$ws.Cells.Item(32, 11).Value = 1846.8334
$ws.Cells.Item(32, 13).Value = -1520.8334
$ws.Cells.Item(39, 8).Value = 1095.6154
$ws.Cells.Item(39, 9).Value = 225.8
$ws.Cells.Item(39, 10).Value = 3995
$ws.Cells.Item(39, 11).Value = 677.4000000000001
$ws.Cells.Item(39, 12).Value = 11985
$ws.Cells.Item(39, 13).Value = -381.4000000000001
$ws.Cells.Item(39, 14).Value = -12577
$ws.Cells.Item(42, 8).Value = 218.07143
$ws.Cells.Item(42, 9).Value = 141.45454
$ws.Cells.Item(42, 10).Value = 499
$ws.Cells.Item(42, 11).Value = 424.36362
$ws.Cells.Item(42, 12).Value = 1497
$ws.Cells.Item(42, 13).Value = -194.36362
$ws.Cells.Item(42, 14).Value = -1957
$ws.Cells.Item(86, 8).Value = 4555.16
$ws.Cells.Item(86, 9).Value = 3374.125
$ws.Cells.Item(86, 10).Value = 5110.9414
$ws.Cells.Item(86, 11).Value = 3374.125
$ws.Cells.Item(86, 12).Value = 5110.9414
$ws.Cells.Item(86, 13).Value = -2251.125
$ws.Cells.Item(86, 14).Value = -7356.9414
$ws.Cells.Item(89, 8).Value = 4555.16
$ws.Cells.Item(89, 9).Value = 3374.125
$ws.Cells.Item(89, 10).Value = 5110.9414
$ws.Cells.Item(89, 11).Value = 16870.625
$ws.Cells.Item(89, 12).Value = 25554.707
$ws.Cells.Item(89, 13).Value = -11254.625
$ws.Cells.Item(89, 14).Value = -36786.70699999999
$ws.Cells.Item(131, 8).Value = 4123.2144
$ws.Cells.Item(131, 9).Value = 1590.625
$ws.Cells.Item(131, 11).Value = 4771.875
$ws.Cells.Item(131, 13).Value = 268.125
$ws.Cells.Item(132, 8).Value = 85411.25
$ws.Cells.Item(132, 9).Value = 91404.42
$ws.Cells.Item(132, 11).Value = 274213.26
$ws.Cells.Item(132, 13).Value = -271683.26
$ws.Cells.Item(138, 8).Value = 2882.1865
$ws.Cells.Item(138, 9).Value = 1395.9286
$ws.Cells.Item(138, 10).Value = 4224.613
$ws.Cells.Item(138, 11).Value = 4187.7858
$ws.Cells.Item(138, 12).Value = 12673.839
$ws.Cells.Item(138, 13).Value = 952.2142000000003
$ws.Cells.Item(138, 14).Value = -22953.839

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(97, 8).Value = 832.6842
$ws.Cells.Item(97, 9).Value = 867.2222
$ws.Cells.Item(97, 10).Value = 211
$ws.Cells.Item(97, 11).Value = 867.2222
$ws.Cells.Item(97, 12).Value = 211
$ws.Cells.Item(97, 13).Value = -371.2222
$ws.Cells.Item(97, 14).Value = -1203
$ws.Cells.Item(132, 8).Value = 1113369.8
$ws.Cells.Item(132, 9).Value = 1335110.4
$ws.Cells.Item(132, 11).Value = 4005331.2
$ws.Cells.Item(132, 13).Value = -4002801.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 935.6875
$ws.Cells.Item(86, 9).Value = 927.625
$ws.Cells.Item(86, 10).Value = 943.75
$ws.Cells.Item(86, 11).Value = 927.625
$ws.Cells.Item(86, 12).Value = 943.75
$ws.Cells.Item(86, 13).Value = 195.375
$ws.Cells.Item(86, 14).Value = -3189.75
$ws.Cells.Item(89, 8).Value = 935.6875
$ws.Cells.Item(89, 9).Value = 927.625
$ws.Cells.Item(89, 10).Value = 943.75
$ws.Cells.Item(89, 11).Value = 4638.125
$ws.Cells.Item(89, 12).Value = 4718.75
$ws.Cells.Item(89, 13).Value = 977.875
$ws.Cells.Item(89, 14).Value = -15950.75
$ws.Cells.Item(94, 8).Value = 1114
$ws.Cells.Item(94, 9).Value = 1173.6316
$ws.Cells.Item(94, 10).Value = 736.3333
$ws.Cells.Item(94, 11).Value = 1173.6316
$ws.Cells.Item(94, 12).Value = 736.3333
$ws.Cells.Item(94, 13).Value = -722.6315999999999
$ws.Cells.Item(94, 14).Value = -1638.3333
$ws.Cells.Item(132, 8).Value = 0
$ws.Cells.Item(132, 10).Value = 0
$ws.Cells.Item(132, 12).Value = 0
$ws.Cells.Item(132, 14).ClearContents()
$ws.Cells.Item(134, 8).Value = 947812.6
$ws.Cells.Item(134, 9).Value = 1036817.4
$ws.Cells.Item(134, 10).Value = 655368.4399999999
$ws.Cells.Item(134, 11).Value = 3110452.2
$ws.Cells.Item(134, 12).Value = 1966105.32
$ws.Cells.Item(134, 13).Value = -3107917.2
$ws.Cells.Item(134, 14).Value = -1971175.32

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 1019
$ws.Cells.Item(22, 9).Value = 749.6667
$ws.Cells.Item(22, 10).Value = 1692.3334
$ws.Cells.Item(22, 11).Value = 749.6667
$ws.Cells.Item(22, 12).Value = 1692.3334
$ws.Cells.Item(22, 13).Value = -399.6667
$ws.Cells.Item(22, 14).Value = -2392.3334
$ws.Cells.Item(86, 8).Value = 161476.16
$ws.Cells.Item(86, 9).Value = 7073.2856
$ws.Cells.Item(86, 10).Value = 341612.84
$ws.Cells.Item(86, 11).Value = 7073.2856
$ws.Cells.Item(86, 12).Value = 341612.84
$ws.Cells.Item(86, 13).Value = -5950.2856
$ws.Cells.Item(86, 14).Value = -343858.84
$ws.Cells.Item(89, 8).Value = 161476.16
$ws.Cells.Item(89, 9).Value = 7073.2856
$ws.Cells.Item(89, 10).Value = 341612.84
$ws.Cells.Item(89, 11).Value = 35366.428
$ws.Cells.Item(89, 12).Value = 1708064.2
$ws.Cells.Item(89, 13).Value = -29750.428
$ws.Cells.Item(89, 14).Value = -1719296.2
$ws.Cells.Item(99, 8).Value = 5658.8
$ws.Cells.Item(99, 10).Value = 3823.75
$ws.Cells.Item(99, 12).Value = 3823.75
$ws.Cells.Item(99, 14).Value = -6819.75
$ws.Cells.Item(126, 8).Value = 5658.8
$ws.Cells.Item(126, 10).Value = 3823.75
$ws.Cells.Item(126, 12).Value = 11471.25
$ws.Cells.Item(126, 14).Value = -16411.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(11, 8).Value = 21.333334
$ws.Cells.Item(11, 9).Value = 10
$ws.Cells.Item(11, 10).Value = 27
$ws.Cells.Item(11, 11).Value = 30
$ws.Cells.Item(11, 12).Value = 81
$ws.Cells.Item(11, 13).Value = 110
$ws.Cells.Item(11, 14).Value = -361

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(97, 8).Value = 1698.3572
$ws.Cells.Item(97, 9).Value = 1874
$ws.Cells.Item(97, 10).Value = 1382.2
$ws.Cells.Item(97, 11).Value = 1874
$ws.Cells.Item(97, 12).Value = 1382.2
$ws.Cells.Item(97, 13).Value = -1378
$ws.Cells.Item(97, 14).Value = -2374.2
$ws.Cells.Item(107, 8).Value = 23667.732
$ws.Cells.Item(107, 9).Value = 70129.39999999999
$ws.Cells.Item(107, 11).Value = 70129.39999999999
$ws.Cells.Item(107, 13).Value = -68209.39999999999
$ws.Cells.Item(126, 8).Value = 1045119.25
$ws.Cells.Item(126, 9).Value = 1669190.8
$ws.Cells.Item(126, 11).Value = 5007572.4
$ws.Cells.Item(126, 13).Value = -5005102.4
$ws.Cells.Item(132, 8).Value = 10645.214
$ws.Cells.Item(132, 9).Value = 10479.625
$ws.Cells.Item(132, 11).Value = 31438.875
$ws.Cells.Item(132, 13).Value = -28908.875

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 4117.5884
$ws.Cells.Item(7, 9).Value = 4000
$ws.Cells.Item(7, 10).Value = 4333.1665
$ws.Cells.Item(7, 11).Value = 4000
$ws.Cells.Item(7, 12).Value = 4333.1665
$ws.Cells.Item(7, 13).Value = -3888
$ws.Cells.Item(7, 14).Value = -4557.1665
$ws.Cells.Item(22, 8).Value = 1458
$ws.Cells.Item(22, 9).Value = 763.3333
$ws.Cells.Item(22, 10).Value = 2500
$ws.Cells.Item(22, 11).Value = 763.3333
$ws.Cells.Item(22, 12).Value = 2500
$ws.Cells.Item(22, 13).Value = -468.3333
$ws.Cells.Item(22, 14).Value = -3090
$ws.Cells.Item(27, 8).Value = 1458
$ws.Cells.Item(27, 9).Value = 763.3333
$ws.Cells.Item(27, 10).Value = 2500
$ws.Cells.Item(27, 11).Value = 763.3333
$ws.Cells.Item(27, 12).Value = 2500
$ws.Cells.Item(27, 13).Value = -656.3333
$ws.Cells.Item(27, 14).Value = -2714
$ws.Cells.Item(40, 8).Value = 2606.4285
$ws.Cells.Item(40, 9).Value = 2621.8
$ws.Cells.Item(40, 11).Value = 2621.8
$ws.Cells.Item(40, 13).Value = -2485.8
$ws.Cells.Item(82, 8).Value = 2483.3333
$ws.Cells.Item(82, 9).Value = 1933.6666
$ws.Cells.Item(82, 10).Value = 3033
$ws.Cells.Item(82, 11).Value = 1933.6666
$ws.Cells.Item(82, 12).Value = 3033
$ws.Cells.Item(82, 13).Value = -1572.6666
$ws.Cells.Item(82, 14).Value = -3755
$ws.Cells.Item(85, 8).Value = 2483.3333
$ws.Cells.Item(85, 9).Value = 1933.6666
$ws.Cells.Item(85, 10).Value = 3033
$ws.Cells.Item(85, 11).Value = 1933.6666
$ws.Cells.Item(85, 12).Value = 3033
$ws.Cells.Item(85, 13).Value = -685.6666
$ws.Cells.Item(85, 14).Value = -5529
$ws.Cells.Item(126, 8).Value = 4117.5884
$ws.Cells.Item(126, 9).Value = 4000
$ws.Cells.Item(126, 10).Value = 4333.1665
$ws.Cells.Item(126, 11).Value = 12000
$ws.Cells.Item(126, 12).Value = 12999.4995
$ws.Cells.Item(126, 13).Value = -9530
$ws.Cells.Item(126, 14).Value = -17939.4995
$ws.Cells.Item(132, 8).Value = 20601.285
$ws.Cells.Item(132, 9).Value = 22993.166
$ws.Cells.Item(132, 11).Value = 68979.49800000001
$ws.Cells.Item(132, 13).Value = -66449.49800000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(119, 8).Value = 140558.4
$ws.Cells.Item(119, 10).Value = 140558.4
$ws.Cells.Item(119, 12).Value = 140558.4
$ws.Cells.Item(119, 14).Value = -150234.4
$ws.Cells.Item(122, 8).Value = 1950.96
$ws.Cells.Item(122, 9).Value = 1663.8823
$ws.Cells.Item(122, 10).Value = 2561
$ws.Cells.Item(122, 11).Value = 4991.6469
$ws.Cells.Item(122, 12).Value = 7683
$ws.Cells.Item(122, 13).Value = -2541.6469
$ws.Cells.Item(122, 14).Value = -12583
$ws.Cells.Item(126, 8).Value = 1819.7894
$ws.Cells.Item(126, 9).Value = 1740.0588
$ws.Cells.Item(126, 10).Value = 2497.5
$ws.Cells.Item(126, 11).Value = 5220.1764
$ws.Cells.Item(126, 12).Value = 7492.5
$ws.Cells.Item(126, 13).Value = -2750.1764
$ws.Cells.Item(126, 14).Value = -12432.5
$ws.Cells.Item(132, 8).Value = 2875999.2
$ws.Cells.Item(132, 9).Value = 3270254.5
$ws.Cells.Item(132, 11).Value = 9810763.5
$ws.Cells.Item(132, 13).Value = -9808233.5
$ws.Cells.Item(136, 8).Value = 6385.689
$ws.Cells.Item(136, 9).Value = 6450.919
$ws.Cells.Item(136, 11).Value = 19352.757
$ws.Cells.Item(136, 13).Value = -16802.757
